$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 94.9812734082397
$ws.Range("F2").Value = 60.64668769716089

$ws.Range("E3").Value = 5.018726591760299
$ws.Range("F3").Value = 91.04477611940298

$ws.Range("E4").Value = 92.91553133514986
$ws.Range("F4").Value = 93.76832844574781

$ws.Range("E5").Value = 7.084468664850137
$ws.Range("F5").Value = 100

$ws.Range("E6").Value = 99.70760233918129
$ws.Range("F6").Value = 21.9941348973607

$ws.Range("E7").Value = 0.2923976608187134
$ws.Range("F7").Value = 100
